$d = $word.ActiveDocument

# 1) Merge the <id>, p053v_1, </id> runs into a single run containing
#    "<id>p053v_1</id>" (formatting of the merged run follows the first
#    run in the matched range, i.e. the Courier New / 7f6000 formatting).
$r1 = $d.Content
$r1.Find.Execute("<id>p053v_1</id>", $true, $false, $false, $false, $false,
                  $true, 1, $false, "<id>p053v_1</id>", 2)

# 2) "for doing this," -> "to do this," while keeping the original two
#    runs (and their distinct formatting) intact: "for"/"to" keeps the
#    rtl-only run, " doing this,"/" do this," keeps the color+rtl run.
$r2 = $d.Content
$r2.Find.Execute("for doing this,", $true, $false, $false, $false, $false,
                  $true, 1, $false, "", 0)
$runFor = $d.Range($r2.Start, $r2.Start + 3)
$runFor.Text = "to"

$r3 = $d.Content
$r3.Find.Execute(" doing this,", $true, $false, $false, $false, $false,
                  $true, 1, $false, "", 0)
$runDoing = $d.Range($r3.Start, $r3.Start + 6)
$runDoing.Text = " do"
